$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the existing data rows (16:17) so the footer
# rows (formerly 22/23) move down to 24/25.
$ws.Range("B18:B19").EntireRow.Insert()

# Duplicate the two worker rows (16:17) into the new rows (18:19), copying
# values + formatting.
$ws.Range("B16:J17").Copy()
$ws.Range("B18:J19").PasteSpecial()

# Center the "Periodo Mora" column for all four worker rows (matches the
# new copy which uses a centered style for column E).
$ws.Range("E16:E19").HorizontalAlignment = -4108

# New period "2509" for the two new worker rows.
$ws.Range("E18").Value = "2509"
$ws.Range("E19").Value = "2509"

# Totals updated to reflect the new period being added.
$ws.Range("E11").Value = 273880
$ws.Range("F13").Value = 2
